$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C4: unprocessed-docs -> tagged-docs
$ws.Range("C4").Value = "uconn-sdp-team11-tagged-docs"

# Row 7
$ws.Range("A7").Value = "Test3PagePDF_Seven_Ways_to_Apply_the_Cyber_Kill_Chain_with_a_Threat_Intelligence_Platform-page-003.pdf"
$ws.Range("B7").Value = "pdf"
$ws.Range("C7").Value = "uconn-sdp-team11-unprocessed-docs"
$ws.Range("D7").Value = "and"
$ws.Range("E7").Value = "cyber"
$ws.Range("F7").Value = "the"

# Row 8
$ws.Range("A8").Value = "Test3PagePDF_Seven_Ways_to_Apply_the_Cyber_Kill_Chain_with_a_Threat_Intelligence_Platform-page-003.pdf"
$ws.Range("B8").Value = "pdf"
$ws.Range("C8").Value = "uconn-sdp-team11-unprocessed-docs"
$ws.Range("D8").Value = "cyber"
$ws.Range("E8").Value = "threat"

# Row 9
$ws.Range("A9").Value = "catch-can-detecting-server-side-request-forgery-attacks-amazon-web-services_13843.pdf"
$ws.Range("B9").Value = "pdf"
$ws.Range("C9").Value = "uconn-sdp-team11-tagged-docs"
$ws.Range("D9").Value = "attack"

# Row 10
$ws.Range("A10").Value = "catch-can-detecting-server-side-request-forgery-attacks-amazon-web-services_13843.pdf"
$ws.Range("B10").Value = "pdf"
$ws.Range("C10").Value = "uconn-sdp-team11-unprocessed-docs"
$ws.Range("D10").Value = "attacks"
$ws.Range("E10").Value = "detection"
$ws.Range("F10").Value = "request"
